$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Activiti)
$ws.Range("B2").Value = 13
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 89
$ws.Range("G2").Value = 1.816697285160237

# Row 3 (che)
$ws.Range("B3").Value = 48
$ws.Range("D3").Value = 53
$ws.Range("E3").Value = 126
$ws.Range("G3").Value = 4.964539007092199

# Row 4 (pinpoint)
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 4
$ws.Range("G4").Value = 0.04908577739599951

# Row 5 (skywalking)
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1159644375724778

# Row 6 (wildfly)
$ws.Range("B6").Value = 202
$ws.Range("D6").Value = 335
$ws.Range("E6").Value = 538
$ws.Range("G6").Value = 3.821294126003268

# Row 7 (storm)
$ws.Range("B7").Value = 3
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 4
$ws.Range("G7").Value = 0.1177163037080636
